# Auto-generated edit script: re-applies the scheduled market-data refresh
# for the Jenova Leve-profit workbook (per-sheet literal value updates).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1683.9697
$ws.Range("I43").Value = 1824.9
$ws.Range("J43").Value = 1467.1538
$ws.Range("K43").Value = 1824.9
$ws.Range("L43").Value = 1467.1538
$ws.Range("M43").Value = -1755.9
$ws.Range("N43").Value = -1605.1538
$ws.Range("H116").Value = 8966.964
$ws.Range("I116").Value = 5377.4736
$ws.Range("J116").Value = 16544.777
$ws.Range("K116").Value = 5377.4736
$ws.Range("L116").Value = 16544.777
$ws.Range("M116").Value = -1935.4736
$ws.Range("N116").Value = -23428.777
$ws.Range("H132").Value = 2925.3
$ws.Range("I132").Value = 3110.4
$ws.Range("K132").Value = 9331.200000000001
$ws.Range("M132").Value = -6801.200000000001
$ws.Range("H138").Value = 5894.0684
$ws.Range("J138").Value = 6828.365
$ws.Range("L138").Value = 20485.095
$ws.Range("N138").Value = -30765.095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4009.6562
$ws.Range("I2").Value = 4351.5186
$ws.Range("K2").Value = 4351.5186
$ws.Range("M2").Value = -4238.5186
$ws.Range("H32").Value = 5077.849
$ws.Range("I32").Value = 4402.64
$ws.Range("K32").Value = 4402.64
$ws.Range("M32").Value = -4115.64
$ws.Range("H45").Value = 1902.3334
$ws.Range("J45").Value = 2703.5
$ws.Range("L45").Value = 2703.5
$ws.Range("N45").Value = -3457.5
$ws.Range("H74").Value = 2793.375
$ws.Range("I74").Value = 3448
$ws.Range("K74").Value = 3448
$ws.Range("M74").Value = -2574
$ws.Range("H77").Value = 2793.375
$ws.Range("I77").Value = 3448
$ws.Range("K77").Value = 17240
$ws.Range("M77").Value = -12872
$ws.Range("H97").Value = 5790.8887
$ws.Range("I97").Value = 4077.25
$ws.Range("J97").Value = 19500
$ws.Range("K97").Value = 4077.25
$ws.Range("L97").Value = 19500
$ws.Range("M97").Value = -3581.25
$ws.Range("N97").Value = -20492
$ws.Range("H116").Value = 4009.6562
$ws.Range("I116").Value = 4351.5186
$ws.Range("K116").Value = 4351.5186
$ws.Range("M116").Value = -2057.5186
$ws.Range("H125").Value = 88940.22
$ws.Range("J125").Value = 88940.22
$ws.Range("L125").Value = 88940.22
$ws.Range("N125").Value = -98780.22
$ws.Range("H132").Value = 4495.273
$ws.Range("J132").Value = 8552.556
$ws.Range("L132").Value = 25657.668
$ws.Range("N132").Value = -30717.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4009.6562
$ws.Range("I3").Value = 4351.5186
$ws.Range("K3").Value = 4351.5186
$ws.Range("M3").Value = -4237.5186
$ws.Range("H94").Value = 2990.6667
$ws.Range("I94").Value = 2958.2856
$ws.Range("J94").Value = 3444
$ws.Range("K94").Value = 2958.2856
$ws.Range("L94").Value = 3444
$ws.Range("M94").Value = -2507.2856
$ws.Range("N94").Value = -4346
$ws.Range("H107").Value = 835530.0600000001
$ws.Range("I107").Value = 1705.3334
$ws.Range("K107").Value = 1705.3334
$ws.Range("M107").Value = 214.6666
$ws.Range("H110").Value = 60000
$ws.Range("J110").Value = 60000
$ws.Range("L110").Value = 60000
$ws.Range("N110").Value = -68180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4332.2
$ws.Range("I16").Value = 4191.96
$ws.Range("J16").Value = 4682.8
$ws.Range("K16").Value = 4191.96
$ws.Range("L16").Value = 4682.8
$ws.Range("M16").Value = -3904.96
$ws.Range("N16").Value = -5256.8
$ws.Range("H99").Value = 6538.077
$ws.Range("I99").Value = 5713.5713
$ws.Range("K99").Value = 5713.5713
$ws.Range("M99").Value = -4215.5713
$ws.Range("H105").Value = 1067
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("H113").Value = 4332.2
$ws.Range("I113").Value = 4191.96
$ws.Range("J113").Value = 4682.8
$ws.Range("K113").Value = 4191.96
$ws.Range("L113").Value = 4682.8
$ws.Range("M113").Value = -2021.96
$ws.Range("N113").Value = -9022.799999999999
$ws.Range("H126").Value = 6538.077
$ws.Range("I126").Value = 5713.5713
$ws.Range("K126").Value = 17140.7139
$ws.Range("M126").Value = -14670.7139
$ws.Range("H134").Value = 373865.28
$ws.Range("I134").Value = 3629.3462
$ws.Range("J134").Value = 10000000
$ws.Range("K134").Value = 10888.0386
$ws.Range("L134").Value = 30000000
$ws.Range("M134").Value = -8353.0386
$ws.Range("N134").Value = -30005070
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 491796.53
$ws.Range("I5").Value = 42695.633
$ws.Range("J5").Value = 3336102.2
$ws.Range("K5").Value = 128086.899
$ws.Range("L5").Value = 10008306.6
$ws.Range("M5").Value = -127974.899
$ws.Range("N5").Value = -10008530.6
$ws.Range("H26").Value = 167233.92
$ws.Range("I26").Value = 222523
$ws.Range("J26").Value = 1366.6666
$ws.Range("K26").Value = 667569
$ws.Range("L26").Value = 4099.9998
$ws.Range("M26").Value = -667281
$ws.Range("N26").Value = -4675.9998
$ws.Range("H39").Value = 11279.7
$ws.Range("I39").Value = 3319
$ws.Range("J39").Value = 13933.267
$ws.Range("K39").Value = 9957
$ws.Range("L39").Value = 41799.801
$ws.Range("M39").Value = -9663
$ws.Range("N39").Value = -42387.801
$ws.Range("H74").Value = 3118.8
$ws.Range("I74").Value = 988
$ws.Range("J74").Value = 6315
$ws.Range("K74").Value = 2964
$ws.Range("L74").Value = 18945
$ws.Range("M74").Value = -1903
$ws.Range("N74").Value = -21067
$ws.Range("H77").Value = 3118.8
$ws.Range("I77").Value = 988
$ws.Range("J77").Value = 6315
$ws.Range("K77").Value = 8892
$ws.Range("L77").Value = 56835
$ws.Range("M77").Value = -3588
$ws.Range("N77").Value = -67443
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("H101").Value = 5314
$ws.Range("J101").Value = 5764.5
$ws.Range("L101").Value = 17293.5
$ws.Range("N101").Value = -22161.5
$ws.Range("H113").Value = 1426321.5
$ws.Range("I113").Value = 3704815
$ws.Range("J113").Value = 2263.0625
$ws.Range("K113").Value = 11114445
$ws.Range("L113").Value = 6789.1875
$ws.Range("M113").Value = -11112275
$ws.Range("N113").Value = -11129.1875
$ws.Range("H132").Value = 503328.53
$ws.Range("J132").Value = 670695.25
$ws.Range("L132").Value = 6036257.25
$ws.Range("N132").Value = -6041317.25
$ws.Range("H135").Value = 491796.53
$ws.Range("I135").Value = 42695.633
$ws.Range("J135").Value = 3336102.2
$ws.Range("K135").Value = 384260.697
$ws.Range("L135").Value = 30024919.8
$ws.Range("M135").Value = -381725.697
$ws.Range("N135").Value = -30029989.8
$ws.Range("H137").Value = 3305.4348
$ws.Range("I137").Value = 2907.9375
$ws.Range("J137").Value = 4214
$ws.Range("K137").Value = 8723.8125
$ws.Range("L137").Value = 12642
$ws.Range("M137").Value = -3623.8125
$ws.Range("N137").Value = -22842
$ws.Range("N80").ClearContents()
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 507980.06
$ws.Range("I113").Value = 1002046.3
$ws.Range("J113").Value = 13913.8
$ws.Range("K113").Value = 1002046.3
$ws.Range("L113").Value = 13913.8
$ws.Range("M113").Value = -999876.3
$ws.Range("N113").Value = -18253.8
$ws.Range("H122").Value = 2719.4849
$ws.Range("I122").Value = 1890.6296
$ws.Range("K122").Value = 5671.8888
$ws.Range("M122").Value = -3221.8888

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1103.8889
$ws.Range("I55").Value = 280.16666
$ws.Range("K55").Value = 280.16666
$ws.Range("M55").Value = -107.16666
$ws.Range("H93").Value = 2234.8
$ws.Range("I93").Value = 1932.7407
$ws.Range("J93").Value = 3254.25
$ws.Range("K93").Value = 1932.7407
$ws.Range("L93").Value = 3254.25
$ws.Range("M93").Value = -684.7407000000001
$ws.Range("N93").Value = -5750.25
$ws.Range("H100").Value = 2330.3333
$ws.Range("I100").Value = 2061.5715
$ws.Range("K100").Value = 2061.5715
$ws.Range("M100").Value = -1520.5715
$ws.Range("H122").Value = 4803.653
$ws.Range("I122").Value = 4510.552
$ws.Range("J122").Value = 5228.65
$ws.Range("K122").Value = 13531.656
$ws.Range("L122").Value = 15685.95
$ws.Range("M122").Value = -11081.656
$ws.Range("N122").Value = -20585.95

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 26318516
$ws.Range("I122").Value = 37039316
$ws.Range("J122").Value = 3820.4546
$ws.Range("K122").Value = 111117948
$ws.Range("L122").Value = 11461.3638
$ws.Range("M122").Value = -111115498
$ws.Range("N122").Value = -16361.3638
$ws.Range("H126").Value = 1048.8823
$ws.Range("I126").Value = 989.4375
$ws.Range("K126").Value = 2968.3125
$ws.Range("M126").Value = -498.3125
